# "Generate Report for Archive" — regenerate the localization-status report:
#  - the "Ready for handoff" status has moved on to "In Translation" for this
#    file, on every sheet that surfaces it (Overview summary + each locale
#    sheet's Status column)
#  - the Status column narrows to fit the new (shorter) text

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.576851254417766

# --- Per-locale sheets: Status column (C) ----------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.576851254417766

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.576851254417766
